# Weekly data update: insert a new price observation row for "Haba"
# (Vega Central Mapocho de Santiago) just before the existing row 167,
# pushing the existing rows 167-237 down to 168-238.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 167 (shifts rows 167..237 -> 168..238)
$ws.Rows(167).Insert()

# Populate the newly inserted row 167 with the new weekly record.
$ws.Cells.Item(167, 1).Value  = 9
$ws.Cells.Item(167, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(167, 3).Value  = "Metropolitana"
$ws.Cells.Item(167, 4).Value  = 44755
$ws.Cells.Item(167, 5).Value  = 13
$ws.Cells.Item(167, 6).Value  = 100112026
$ws.Cells.Item(167, 7).Value  = "Haba"
$ws.Cells.Item(167, 8).Value  = "Sin especificar"
$ws.Cells.Item(167, 9).Value  = "Primera"
$ws.Cells.Item(167, 10).Value = 52
$ws.Cells.Item(167, 11).Value = 16000
$ws.Cells.Item(167, 12).Value = 18000
$ws.Cells.Item(167, 13).Value = 17000
$ws.Cells.Item(167, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(167, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(167, 16).Value = 680
$ws.Cells.Item(167, 17).Value = 25
$ws.Cells.Item(167, 18).Value = "Hortaliza"

# Keep the date cell formatted the same way as its neighbours.
$ws.Cells.Item(167, 4).NumberFormat = $ws.Cells.Item(168, 4).NumberFormat
